$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 32: fill in date/hours, add new activity string ---
$ws.Range("A32").NumberFormat = "mm/dd/yy"
$ws.Range("A32").Value = 44252
$ws.Range("B32").Value = 8
$ws.Range("D32").Value = "Debugging mit Henneke bezüglich Senden der Decoderliste vom RedPitaya zum Client"

# --- Row 33 ---
$ws.Range("A33").NumberFormat = "mm/dd/yy"
$ws.Range("A33").Value = 44253
$ws.Range("B33").Value = 8

# --- Row 34 ---
$ws.Range("A34").NumberFormat = "mm/dd/yy"
$ws.Range("A34").Value = 44254
$ws.Range("B34").Value = 8

# --- Row 35 ---
$ws.Range("A35").NumberFormat = "mm/dd/yy"
$ws.Range("A35").Value = 44255
$ws.Range("B35").Value = 8
$ws.Rows(35).RowHeight = 13.8

# --- Row 36 (new) ---
$ws.Range("A36").NumberFormat = "mm/dd/yy"
$ws.Range("A36").Value = 44230
$ws.Range("B36").Value = 8
$ws.Range("C36").Formula = "=C35+B36"
$ws.Rows(36).RowHeight = 13.8

# --- Row 37 (new) ---
$ws.Range("A37").NumberFormat = "mm/dd/yy"
$ws.Range("A37").Value = 44231
$ws.Range("B37").Value = 8
$ws.Range("C37").Formula = "=C36+B37"
$ws.Rows(37).RowHeight = 13.8

# --- Update selection to match the author's final cursor position ---
$null = $ws.Range("D33").Select()
